# driver-wait untill ele visible
# Adds the "Nom_ad1 / Nom_ad2 / Nom_pincode" nominee-address columns
# (AO:AQ) to row 1 (headers) and row 2 (data), and switches the existing
# R2 pincode cell from a literal number to the shared "Nom_pincode" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- build the header font (Arial Unicode MS, 10pt) via a throwaway named
# style so the engine materializes the font in one shared-font hop instead
# of the multi-hop "de-schemed clone" path that Range.Font.Name takes ---
$tmpStyle = $wb.Styles.Add("NomHeaderTmp")
$tmpStyle.Font.Name = "Arial Unicode MS"
$ws.Range("AO1").Style = "NomHeaderTmp"
$tmpStyle.Delete()
$ws.Range("AO1").Font.Size = 10
$ws.Range("AO1").VerticalAlignment = -4108

# --- write the new cells in the same order the original edit's shared
# strings table was built in (Nom_ad1, Vadapalani, Nom_ad2, street, Nom_pincode, 777755553214) ---
$ws.Range("AO1").Value = "Nom_ad1"
$ws.Range("AP2").Value = "Vadapalani"
$ws.Range("AP1").Value = "Nom_ad2"
$ws.Range("AO2").Value = "Gangai amman kovil street"
$ws.Range("AQ1").Value = "Nom_pincode"
$ws.Range("R2").Value = "777755553214"

# AQ2 keeps a genuine numeric value (not shared-string text) even though
# the column's default style is text-formatted (numFmtId 49): reset to the
# Normal style first, write the number, then restore the text format code
# so the stored <v> stays numeric.
$ws.Range("AQ2").Style = "Normal"
$ws.Range("AQ2").Value = 600026
$ws.Range("AQ2").NumberFormat = "@"

# keep the selection where the author's last edit left it
$ws.Range("S2").Select() | Out-Null
